# Update the presenter-notes on slide 1 ("Note from Chuck ...") so that the
# note is given in Greek, with "Chuck" kept in Latin script, matching the
# author's edit.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$notes = $s.NotesPage

# The note lives in the body placeholder of the notes page (first shape).
$shp = $notes.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

$run1 = "Σημείωση από τον "
$run2 = " Chuck"
$run3 = ". Εάν χρησιμοποιείτε αυτό το υλικό, μπορείτε να αφαιρέσετε το λογότυπο UM και να το αντικαταστήσετε με το δικό σας, αλλά διατηρήστε το λογότυπο CC-BY στην πρώτη σελίδα καθώς την/τις σελίδα/"
$run4 = "ες"
$run5 = " "
$run6 = "αναγνώρισης."
$run7 = "."

$tr.Text = $run1 + $run2 + $run3 + $run4 + $run5 + $run6 + $run7
